$wb = $excel.ActiveWorkbook

# Rename sheets: hyphens to underscores
$wb.Worksheets.Item("experiment-description").Name = "experiment_description"
$wb.Worksheets.Item("experiment-specification").Name = "experiment_specification"
$wb.Worksheets.Item("run-description").Name = "run_description"
$wb.Worksheets.Item("run-specification").Name = "run_specification"
$wb.Worksheets.Item("heuristics-description").Name = "heuristics_description"

# Activate the heuristics_description sheet (5th sheet) so it becomes the
# selected/active tab, and update its selected cell.
$ws5 = $wb.Worksheets.Item("heuristics_description")
$ws5.Activate()
$ws5.Range("J23").Select()
